$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new status-label strings in the exact order they should land in
# the shared-strings table (first-use order), then fill in the rest of the
# column E cells.
$ws.Range("E2").Value = "Alias Changed for Network"
$ws.Range("E21").Value = "Network Removed from Database"
$ws.Range("E23").Value = "New Network Added to Database"
$ws.Range("E38").Value = "Network Moved from One Add-On Package to Another Add-On Package"
$ws.Range("E7").Value = "Network Added to Base Service"
$ws.Range("E6").Value = "Network Added to Add-On Package"
$ws.Range("E37").Value = "Network Removed from Base Service"
$ws.Range("E39").Value = "Name of Add-On Package Changed"

$ws.Range("E3").Value = "Alias Changed for Network"
$ws.Range("E4").Value = "Alias Changed for Network"
$ws.Range("E5").Value = "Alias Changed for Network"
$ws.Range("E8").Value = "Network Added to Add-On Package"
$ws.Range("E9").Value = "Network Added to Base Service"
$ws.Range("E10").Value = "Network Added to Add-On Package"
$ws.Range("E11").Value = "Network Added to Add-On Package"
$ws.Range("E12").Value = "Network Added to Add-On Package"
$ws.Range("E13").Value = "Network Added to Add-On Package"
$ws.Range("E14").Value = "Network Added to Base Service"
$ws.Range("E15").Value = "Network Added to Base Service"
$ws.Range("E16").Value = "Network Added to Base Service"
$ws.Range("E17").Value = "Network Added to Add-On Package"
$ws.Range("E18").Value = "Network Added to Base Service"
$ws.Range("E19").Value = "Network Added to Base Service"
$ws.Range("E20").Value = "Network Added to Add-On Package"
$ws.Range("E22").Value = "Network Removed from Database"
$ws.Range("E24").Value = "Network Removed from Database"
$ws.Range("E25").Value = "Network Removed from Database"
$ws.Range("E26").Value = "New Network Added to Database"
$ws.Range("E27").Value = "Network Removed from Database"
$ws.Range("E28").Value = "New Network Added to Database"
$ws.Range("E29").Value = "New Network Added to Database"
$ws.Range("E30").Value = "Network Removed from Database"
$ws.Range("E31").Value = "Network Removed from Database"
$ws.Range("E32").Value = "New Network Added to Database"
$ws.Range("E33").Value = "Network Removed from Database"
$ws.Range("E34").Value = "New Network Added to Database"
$ws.Range("E35").Value = "New Network Added to Database"
$ws.Range("E36").Value = "New Network Added to Database"
$ws.Range("E40").Value = "Name of Add-On Package Changed"
$ws.Range("E41").Value = "Network Removed from Base Service"
$ws.Range("E42").Value = "Name of Add-On Package Changed"

# Matches the saved selection state captured in the workbook (E39:E42 range
# selected, active cell E39).
$ws.Range("E39:E42").Select()
